$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022399915397922
$ws.Range("D2").Value = 1.02763747109014
$ws.Range("E2").Value = 1.023155570628548
$ws.Range("F2").Value = 1.035986868834263
$ws.Range("I2").Value = 1.032199466818722
$ws.Range("J2").Value = 1.027585967381304
$ws.Range("K2").Value = 1.030456618849389
$ws.Range("L2").Value = 1.025987839658773
$ws.Range("M2").Value = 1.038781901741791
$ws.Range("N2").Value = 1.013264242149561
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023135545098344
$ws.Range("D3").Value = 1.028179636389429
$ws.Range("E3").Value = 1.023774320408744
$ws.Range("F3").Value = 1.038101791240905
$ws.Range("I3").Value = 1.032387479197491
$ws.Range("J3").Value = 1.027960741005537
$ws.Range("K3").Value = 1.030807269282013
$ws.Range("L3").Value = 1.026413921480994
$ws.Range("M3").Value = 1.040702870266985
$ws.Range("N3").Value = 1.013389081851469
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023612029969167
$ws.Range("D4").Value = 1.028530759865954
$ws.Range("E4").Value = 1.024175503470641
$ws.Range("F4").Value = 1.039465354565955
$ws.Range("I4").Value = 1.032507873646105
$ws.Range("J4").Value = 1.028203042437432
$ws.Range("K4").Value = 1.031033761426051
$ws.Range("L4").Value = 1.02668974355078
$ws.Range("M4").Value = 1.041940587210116
$ws.Range("N4").Value = 1.013469770746543
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023812458854264
$ws.Range("D5").Value = 1.02867844470796
$ws.Range("E5").Value = 1.024344353668992
$ws.Range("F5").Value = 1.040037444298577
$ws.Range("I5").Value = 1.032558186084092
$ws.Range("J5").Value = 1.028304857133347
$ws.Range("K5").Value = 1.031128882229427
$ws.Range("L5").Value = 1.026805726865936
$ws.Range("M5").Value = 1.042459686440195
$ws.Range("N5").Value = 1.013503670457349
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023846118438106
$ws.Range("D6").Value = 1.028703245872756
$ws.Range("E6").Value = 1.02437271563281
$ws.Range("F6").Value = 1.040133433936862
$ws.Range("I6").Value = 1.032566616112232
$ws.Range("J6").Value = 1.028321949392915
$ws.Range("K6").Value = 1.031144847765628
$ws.Range("L6").Value = 1.02682520258005
$ws.Range("M6").Value = 1.04254677371548
$ws.Range("N6").Value = 1.013509361078443
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023614707660604
$ws.Range("D7").Value = 1.028532732954075
$ws.Range("E7").Value = 1.024177758899091
$ws.Range("F7").Value = 1.039473003349779
$ws.Range("I7").Value = 1.032508547106693
$ws.Range("J7").Value = 1.028204403082801
$ws.Range("K7").Value = 1.031035032814352
$ws.Range("L7").Value = 1.026691293217495
$ws.Range("M7").Value = 1.041947528267501
$ws.Range("N7").Value = 1.013470223802425
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022648424209795
$ws.Range("D8").Value = 1.027820634506422
$ws.Range("E8").Value = 1.023364511452125
$ws.Range("F8").Value = 1.036702655506081
$ws.Range("I8").Value = 1.032263268284127
$ws.Range("J8").Value = 1.027712665488712
$ws.Range("K8").Value = 1.030575205940606
$ws.Range("L8").Value = 1.026131811032427
$ws.Range("M8").Value = 1.039432211049672
$ws.Range("N8").Value = 1.013306450994458
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020949456014173
$ws.Range("D9").Value = 1.026568205473672
$ws.Range("E9").Value = 1.021937726452428
$ws.Range("F9").Value = 1.031781871526001
$ws.Range("I9").Value = 1.031821353692609
$ws.Range("J9").Value = 1.026844630696007
$ws.Range("K9").Value = 1.029761861105709
$ws.Range("L9").Value = 1.025146870054826
$ws.Range("M9").Value = 1.034958271059596
$ws.Range("N9").Value = 1.013017174958473
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019819388826343
$ws.Range("D10").Value = 1.025734895614699
$ws.Range("E10").Value = 1.020990815605172
$ws.Range("F10").Value = 1.028473247766501
$ws.Range("I10").Value = 1.031520169071237
$ws.Range("J10").Value = 1.026264934213037
$ws.Range("K10").Value = 1.029217575630128
$ws.Range("L10").Value = 1.024490911538607
$ws.Range("M10").Value = 1.031945955861773
$ws.Range("N10").Value = 1.012823871593205
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019330680613148
$ws.Range("D11").Value = 1.025374461228677
$ws.Range("E11").Value = 1.020581821929762
$ws.Range("F11").Value = 1.027033501093418
$ws.Range("I11").Value = 1.031388181130102
$ws.Range("J11").Value = 1.026013684546025
$ws.Range("K11").Value = 1.028981407613776
$ws.Range("L11").Value = 1.024207040392957
$ws.Range("M11").Value = 1.030634166371325
$ws.Range("N11").Value = 1.01274006353383
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019149246413652
$ws.Range("D12").Value = 1.025240639692226
$ws.Range("E12").Value = 1.020430058572179
$ws.Range("F12").Value = 1.02649761411653
$ws.Range("I12").Value = 1.031338917483726
$ws.Range("J12").Value = 1.025920323852178
$ws.Range("K12").Value = 1.028893610872203
$ws.Range("L12").Value = 1.024101623089896
$ws.Range("M12").Value = 1.030145758833422
$ws.Range("N12").Value = 1.012708917624571
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019188160406992
$ws.Range("D13").Value = 1.025269342148602
$ws.Range("E13").Value = 1.020462605291121
$ws.Range("F13").Value = 1.02661261409143
$ws.Range("I13").Value = 1.031349495466452
$ws.Range("J13").Value = 1.025940351634916
$ws.Range("K13").Value = 1.028912446893961
$ws.Range("L13").Value = 1.024124234315579
$ws.Range("M13").Value = 1.030250576509708
$ws.Range("N13").Value = 1.0127155992438
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01931568128723
$ws.Range("D14").Value = 1.025363398274548
$ws.Range("E14").Value = 1.020569273952848
$ws.Range("F14").Value = 1.02698922717632
$ws.Range("I14").Value = 1.031384113832191
$ws.Range("J14").Value = 1.026005968042237
$ws.Range("K14").Value = 1.028974151805582
$ws.Range("L14").Value = 1.024198326051024
$ws.Range("M14").Value = 1.030593818115883
$ws.Range("N14").Value = 1.012737489326187
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019394263564785
$ws.Range("D15").Value = 1.025421357305402
$ws.Range("E15").Value = 1.020635016621553
$ws.Range("F15").Value = 1.027221123963387
$ws.Range("I15").Value = 1.031405411851618
$ws.Range("J15").Value = 1.026046391795733
$ws.Range("K15").Value = 1.02901216050864
$ws.Range("L15").Value = 1.024243979754239
$ws.Range("M15").Value = 1.030805147141843
$ws.Range("N15").Value = 1.012750974429836
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01985183588521
$ws.Range("D16").Value = 1.025758824817307
$ws.Range("E16").Value = 1.021017980858931
$ws.Range("F16").Value = 1.028568646587739
$ws.Range("I16").Value = 1.031528895439642
$ws.Range("J16").Value = 1.026281603850201
$ws.Range("K16").Value = 1.029233239026547
$ws.Range("L16").Value = 1.024509754618546
$ws.Range("M16").Value = 1.03203285549951
$ws.Range("N16").Value = 1.012829431426623
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.020139025095052
$ws.Range("D17").Value = 1.025970615268788
$ws.Range("E17").Value = 1.021258479413677
$ws.Range("F17").Value = 1.029411988548154
$ws.Range("I17").Value = 1.031605931509876
$ws.Range("J17").Value = 1.026429082828776
$ws.Range("K17").Value = 1.02937178488762
$ws.Range("L17").Value = 1.024676512267459
$ws.Range("M17").Value = 1.032800950141429
$ws.Range("N17").Value = 1.01287861702472
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020306597304905
$ws.Range("D18").Value = 1.026094187059995
$ws.Range("E18").Value = 1.021398856973637
$ws.Range("F18").Value = 1.029903213927805
$ws.Range("I18").Value = 1.03165071365105
$ws.Range("J18").Value = 1.026515081913104
$ws.Range("K18").Value = 1.029452549181968
$ws.Range("L18").Value = 1.024773794895357
$ws.Range("M18").Value = 1.033248251443056
$ws.Range("N18").Value = 1.012907295875498
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020363745179853
$ws.Range("D19").Value = 1.026136328274412
$ws.Range("E19").Value = 1.021446738821352
$ws.Range("F19").Value = 1.030070594527
$ws.Range("I19").Value = 1.031665957507092
$ws.Range("J19").Value = 1.026544401484116
$ws.Range("K19").Value = 1.029480079718608
$ws.Range("L19").Value = 1.024806968405122
$ws.Range("M19").Value = 1.033400649294639
$ws.Range("N19").Value = 1.012917072875444
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.020108206239458
$ws.Range("D20").Value = 1.025947888222434
$ws.Range("E20").Value = 1.021232665965222
$ws.Range("F20").Value = 1.029321576706755
$ws.Range("I20").Value = 1.031597681967236
$ws.Range("J20").Value = 1.02641326209778
$ws.Range("K20").Value = 1.029356925106068
$ws.Range("L20").Value = 1.024658619120699
$ws.Range("M20").Value = 1.032718615081703
$ws.Range("N20").Value = 1.012873340937583
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019278126970225
$ws.Range("D21").Value = 1.025335699433055
$ws.Range("E21").Value = 1.020537858400154
$ws.Range("F21").Value = 1.026878354690708
$ws.Range("I21").Value = 1.031373926150208
$ws.Range("J21").Value = 1.025986646620369
$ws.Range("K21").Value = 1.028955983273957
$ws.Range("L21").Value = 1.024176507184913
$ws.Range("M21").Value = 1.030492774046837
$ws.Range("N21").Value = 1.012731043679486
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018756766258462
$ws.Range("D22").Value = 1.024951138555325
$ws.Range("E22").Value = 1.020101903225563
$ws.Range("F22").Value = 1.025335816145898
$ws.Range("I22").Value = 1.031231867607256
$ws.Range("J22").Value = 1.025718211952941
$ws.Range("K22").Value = 1.028703470406141
$ws.Range("L22").Value = 1.023873529641923
$ws.Range("M22").Value = 1.029086626095992
$ws.Range("N22").Value = 1.012641484015555
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019033097666235
$ws.Range("D23").Value = 1.025154968511969
$ws.Range("E23").Value = 1.020332925833487
$ws.Range("F23").Value = 1.0261541625699
$ws.Range("I23").Value = 1.031307306186057
$ws.Range("J23").Value = 1.02586053355755
$ws.Range("K23").Value = 1.028837372524888
$ws.Range("L23").Value = 1.024034129809479
$ws.Range("M23").Value = 1.029832695504116
$ws.Range("N23").Value = 1.012688969942857
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.020122131770021
$ws.Range("D24").Value = 1.025958157479322
$ws.Range("E24").Value = 1.021244329646503
$ws.Range("F24").Value = 1.029362432031334
$ws.Range("I24").Value = 1.031601410049264
$ws.Range("J24").Value = 1.026420410876153
$ws.Range("K24").Value = 1.029363639747668
$ws.Range("L24").Value = 1.02466670421525
$ws.Range("M24").Value = 1.03275582096019
$ws.Range("N24").Value = 1.012875725006092
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021388229848822
$ws.Range("D25").Value = 1.026891701899785
$ws.Range("E25").Value = 1.022305835374322
$ws.Range("F25").Value = 1.033058829420263
$ws.Range("I25").Value = 1.0319367541836
$ws.Range("J25").Value = 1.027069217488174
$ws.Range("K25").Value = 1.029972493154349
$ws.Range("L25").Value = 1.025401385750917
$ws.Range("M25").Value = 1.036119998416571
$ws.Range("N25").Value = 1.013092040208938
